# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (copied from the existing "2022-Q2"
# sheet so it keeps the same layout/formatting), fills it with the new
# quarter's fund data, and updates the "总计" (totals) summary sheet with
# a new row for 2022-Q3 while keeping the existing 2022-Q2 / 2022-Q1 rows.

$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell as TEXT without Excel auto-converting
# numeric-looking strings (e.g. "4.99", or fund codes with leading zeros
# like "003721") into numbers, and without leaving behind a stray
# NumberFormat/style change on the target cell. We stage the text in a
# far-away scratch cell (forced to Text format), copy it, and paste only
# the VALUE into the destination - the destination keeps its own style.
function Set-TextValue {
    param($sheet, $addr, $val)
    $scratch = $sheet.Range("Z100")
    $scratch.NumberFormat = "@"
    $scratch.Value = $val
    $scratch.Copy()
    $sheet.Range($addr).PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = $false
    $scratch.Clear()
}

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by duplicating "2022-Q2" (index 2)
#    and placing the copy immediately before it. This keeps identical
#    column widths / styles / borders to the existing quarter sheets.
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item(2)
$q2Sheet.Copy($q2Sheet)

$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

Set-TextValue $q3Sheet "B2" "161128"
Set-TextValue $q3Sheet "C2" "易方达标普信息科技指数（QDII-LOF）人民币"
Set-TextValue $q3Sheet "D2" "4.99"
Set-TextValue $q3Sheet "E2" "91.96"
Set-TextValue $q3Sheet "F2" "1.87"
Set-TextValue $q3Sheet "G2" "0.0933"
$q3Sheet.Range("H2").Value = 8

Set-TextValue $q3Sheet "B3" "012868"
Set-TextValue $q3Sheet "C3" "易方达标普信息科技指数（QDII-LOF）人民币 C"
Set-TextValue $q3Sheet "D3" "4.99"
Set-TextValue $q3Sheet "E3" "91.96"
Set-TextValue $q3Sheet "F3" "1.87"
Set-TextValue $q3Sheet "G3" "0.0933"
$q3Sheet.Range("H3").Value = 8

Set-TextValue $q3Sheet "B4" "003721"
Set-TextValue $q3Sheet "C4" "易方达标普信息科技指数（QDII-LOF）美元A"
Set-TextValue $q3Sheet "D4" "4.84"
Set-TextValue $q3Sheet "E4" "91.96"
Set-TextValue $q3Sheet "F4" "1.87"
Set-TextValue $q3Sheet "G4" "0.0905"
$q3Sheet.Range("H4").Value = 8

Set-TextValue $q3Sheet "B5" "012869"
Set-TextValue $q3Sheet "C5" "易方达标普信息科技指数（QDII-LOF）美元 C"
Set-TextValue $q3Sheet "D5" "0.15"
Set-TextValue $q3Sheet "E5" "91.96"
Set-TextValue $q3Sheet "F5" "1.87"
Set-TextValue $q3Sheet "G5" "0.0028"
$q3Sheet.Range("H5").Value = 8

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row for 2022-Q3 above
#    the existing 2022-Q2 row, and adjust the 2022-Q2 holdings figure.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Insert() drags a blended style (border/alignment) down from the header
# row into the freshly inserted row. Clear that and copy the plain
# formatting from the row immediately below (the original data row,
# now shifted to row 3) so the new row matches the other data rows.
$totalSheet.Range("A2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$totalSheet.Range("A2").Value = 0
Set-TextValue $totalSheet "B2" "2022-Q3"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.28

$totalSheet.Range("A3").Value = 1
Set-TextValue $totalSheet "B3" "2022-Q2"
$totalSheet.Range("C3").Value = 4
$totalSheet.Range("D3").Value = 0.2

$totalSheet.Range("A4").Value = 2
Set-TextValue $totalSheet "B4" "2022-Q1"
$totalSheet.Range("C4").Value = 2
$totalSheet.Range("D4").Value = 0.23
